$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people): 2.7 -> 2.65
$ws.Range("C13").Value = "'2.65"
$ws.Range("C13").Style = "Normal"

# Value added to the economy (% of total): Micro/SMEs/MSMEs
$ws.Range("B18").Value = "'34.73"
$ws.Range("B18").Style = "Normal"

$ws.Range("C18").Value = "'23.21"
$ws.Range("C18").Style = "Normal"

$ws.Range("D18").Value = "'57.94"
$ws.Range("D18").Style = "Normal"
